# Update the "Development Process" slide (slide 5) by setting the text
# of its content placeholder shape (id=3, "Content Placeholder 2") to
# "Requirement Analysis".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Requirement Analysis"
$shape.TextFrame.TextRange.LanguageID = "en-CA"
